$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 0.7037333250045776
$ws.Range("B1").Value = 0.8482045531272888
$ws.Range("C1").Value = 1.147632598876953
$ws.Range("D1").Value = 3.269025564193726
$ws.Range("E1").Value = 4.278785228729248
